# Insert four new columns of data into row 1 of the active sheet: the
# existing "id_centro" header (previously in I1) moves to M1, and the
# newly inserted I1:L1 cells get the new header names in between.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "id_centro" header text/format by moving it to M1
# before overwriting I1 with the new header.
$ws.Range("M1").Value = "id_centro"
$ws.Range("I1").Value = "cr_total_dias_ingreso"
$ws.Range("J1").Value = "primera_conexion_crea"
$ws.Range("K1").Value = "dias_de_conexion_dispositivo"
$ws.Range("L1").Value = "primera_conexion_dispositivo"

# Match the header style (bold font, thin border, centered/top-aligned)
# used by the rest of row 1 by copying the format from an existing
# header cell onto the newly populated ones.
$ws.Range("H1").Copy()
$ws.Range("I1:M1").PasteSpecial(-4122)
